$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (rows 2-11) so stale rows/shared strings are dropped.
$ws.Range("A2:I11").ClearContents()

# New survey data: all rows now belong to a single refined userId,
# with planetNo/attemptNo 1/1, 1/2, 2/1, 3/1, 3/2.
$userId = "09ee0a844a2345f393f2502ef4514d5b"

$data = @(
    @($userId, 1, 1, 27, 63, 0, 29, 61, 0),
    @($userId, 1, 2, 29, 61, 1, 29, 61, 1),
    @($userId, 2, 1, 47, 14, 1, 47, 14, 1),
    @($userId, 3, 1,  0, 14, 0, 35, -21, 0),
    @($userId, 3, 2, 12, 78, 1, 35,  55, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}
